$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "RDJ" -> "rdj" (the actor name cell lost its uppercase styling)
$ws.Range("B2").Value = "rdj"

# move the active selection to B2, matching the saved view state
[void]$ws.Range("B2").Select()
